# The second paragraph of the document holds a M2Doc query written as a
# Word field:  { m:'doc.html'.fromHTMLURI() }
#
# Originally that query is stored as a real Word field: a fldChar "begin",
# a run of w:instrText pieces carrying the query text (split across several
# runs, with a bookmark sitting in the middle of them), and a fldChar "end".
#
# The parser was updated (TokenIteratorFieldRewriterSplit) to instead read
# the query straight out of plain w:t runs, so this edit "flattens" the
# field: the fldChar begin/end runs disappear, the leading/trailing space
# instrText runs are folded into literal "{" / "}" text runs, and every
# remaining w:instrText run becomes a plain w:t run with the exact same
# text - the bookmark in the middle is left untouched.
#
# We rebuild the whole paragraph's XML (preserving the existing paragraph
# and run rsid attributes) and push it in with Range.InsertXML, which
# replaces the exact range's contents in one shot.

$d = $word.ActiveDocument

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetParagraph = $p
        break
    }
}
if ($targetParagraph -eq $null) {
    # Fallback: the query lives in the second paragraph of this document.
    $targetParagraph = $d.Paragraphs(2)
}

$newParagraphXml = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' + `
    '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r>' + `
    '<w:r w:rsidR="002033E1"><w:t>:</w:t></w:r>' + `
    '<w:r w:rsidR="008B76C9"><w:t>''</w:t></w:r>' + `
    '<w:r w:rsidR="00E806A4"><w:t>doc.html</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r w:rsidR="008B76C9"><w:t>''.fromHTMLURI()</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
    '</w:p>'

[void]$targetParagraph.Range.InsertXML($newParagraphXml)
